$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1719298245614035
$ws.Range("C2").Value = 0.6105263157894737
$ws.Range("P2").Value = 0.1263157894736842
$ws.Range("S2").Value = 0.0912280701754386
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.0223463687150838
$ws.Range("J3").Value = 0.01675977653631285
$ws.Range("P3").Value = 0.770949720670391
$ws.Range("S3").Value = 0.1843575418994413
$ws.Range("O4").Value = 0.02777777777777778
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("P5").Value = 0.4
$ws.Range("S5").Value = 0.6
$ws.Range("B6").Value = 0.04624277456647399
$ws.Range("D6").Value = 0.005780346820809248
$ws.Range("F6").Value = 0.06936416184971098
$ws.Range("J6").Value = 0.2543352601156069
$ws.Range("O6").Value = 0.005780346820809248
$ws.Range("Q6").Value = 0.1040462427745665
$ws.Range("R6").Value = 0.115606936416185
$ws.Range("S6").Value = 0.3988439306358382
$ws.Range("B7").Value = 0.1271676300578035
$ws.Range("D7").Value = 0.005780346820809248
$ws.Range("F7").Value = 0.03468208092485549
$ws.Range("J7").Value = 0.09248554913294797
$ws.Range("O7").Value = 0.02312138728323699
$ws.Range("Q7").Value = 0.1213872832369942
$ws.Range("R7").Value = 0.1098265895953757
$ws.Range("S7").Value = 0.4855491329479769
$ws.Range("B8").Value = 0.1353211009174312
$ws.Range("D8").Value = 0.01834862385321101
$ws.Range("E8").Value = 0.006880733944954129
$ws.Range("F8").Value = 0.04587155963302753
$ws.Range("J8").Value = 0.09174311926605505
$ws.Range("O8").Value = 0.02752293577981652
$ws.Range("Q8").Value = 0.1422018348623853
$ws.Range("R8").Value = 0.1077981651376147
$ws.Range("S8").Value = 0.4243119266055046
$ws.Range("B9").Value = 0.09202453987730061
$ws.Range("D9").Value = 0.04294478527607362
$ws.Range("E9").Value = 0.006134969325153374
$ws.Range("F9").Value = 0.04294478527607362
$ws.Range("J9").Value = 0.07975460122699386
$ws.Range("O9").Value = 0.01840490797546012
$ws.Range("Q9").Value = 0.1165644171779141
$ws.Range("R9").Value = 0.08588957055214724
$ws.Range("S9").Value = 0.5153374233128835
$ws.Range("B10").Value = 0.1359649122807018
$ws.Range("D10").Value = 0.01973684210526316
$ws.Range("E10").Value = 0.002192982456140351
$ws.Range("F10").Value = 0.0668859649122807
$ws.Range("J10").Value = 0.08662280701754387
$ws.Range("O10").Value = 0.02850877192982456
$ws.Range("Q10").Value = 0.1589912280701754
$ws.Range("R10").Value = 0.1129385964912281
$ws.Range("S10").Value = 0.3881578947368421
$ws.Range("G11").Value = 0.1062992125984252
$ws.Range("J11").Value = 0.09055118110236221
$ws.Range("K11").Value = 0.1535433070866142
$ws.Range("L11").Value = 0.6181102362204725
$ws.Range("S11").Value = 0.03149606299212598
$ws.Range("G12").Value = 0.6848484848484848
$ws.Range("J12").Value = 0.2121212121212121
$ws.Range("L12").Value = 0.01212121212121212
$ws.Range("G13").Value = 0.7254901960784313
$ws.Range("J13").Value = 0.2156862745098039
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.02150537634408602
$ws.Range("H15").Value = 0.1720430107526882
$ws.Range("I15").Value = 0.06989247311827956
$ws.Range("J15").Value = 0.3118279569892473
$ws.Range("K15").Value = 0.03763440860215054
$ws.Range("M15").Value = 0.01612903225806452
$ws.Range("O15").Value = 0.04838709677419355
$ws.Range("S15").Value = 0.3225806451612903
$ws.Range("F16").Value = 0.02512562814070352
$ws.Range("H16").Value = 0.1557788944723618
$ws.Range("I16").Value = 0.1055276381909548
$ws.Range("J16").Value = 0.3417085427135678
$ws.Range("K16").Value = 0.06532663316582915
$ws.Range("M16").Value = 0.02010050251256281
$ws.Range("N16").Value = 0.01005025125628141
$ws.Range("O16").Value = 0.06532663316582915
$ws.Range("S16").Value = 0.2110552763819095
$ws.Range("F17").Value = 0.03065134099616858
$ws.Range("H17").Value = 0.2030651340996169
$ws.Range("I17").Value = 0.08045977011494253
$ws.Range("J17").Value = 0.3563218390804598
$ws.Range("K17").Value = 0.103448275862069
$ws.Range("M17").Value = 0.01915708812260536
$ws.Range("O17").Value = 0.04597701149425287
$ws.Range("S17").Value = 0.1609195402298851
$ws.Range("F18").Value = 0.03883495145631068
$ws.Range("H18").Value = 0.1844660194174757
$ws.Range("I18").Value = 0.09223300970873786
$ws.Range("J18").Value = 0.3300970873786408
$ws.Range("K18").Value = 0.1359223300970874
$ws.Range("M18").Value = 0.01941747572815534
$ws.Range("O18").Value = 0.05825242718446602
$ws.Range("S18").Value = 0.1407766990291262
$ws.Range("F19").Value = 0.01862828111769687
$ws.Range("H19").Value = 0.2413209144792549
$ws.Range("I19").Value = 0.07620660457239628
$ws.Range("J19").Value = 0.3149872988992379
$ws.Range("K19").Value = 0.11346316680779
$ws.Range("M19").Value = 0.03048264182895851
$ws.Range("O19").Value = 0.05927180355630821
$ws.Range("S19").Value = 0.1456392887383573
